$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "Approx." column (D).
# Old D (Approx.) -> E, old E (Comment) -> F.
$ws.Range("D1").EntireColumn.Insert()

# New column header + width
$ws.Range("D1").Value = "Provided dim.[mm]"
$ws.Range("D1").EntireColumn.ColumnWidth = 23.83

# Fill in "Provided dim.[mm]" values for each material row
$ws.Range("D2").Value = "25, 50, 100"
$ws.Range("D3").Value = "100, 120, 140, 160"
$ws.Range("D4").Value = "120, 140, 160, 180, 200, 220, 240, 260, 280, 300"
$ws.Range("D5").Value = "60, 80, 100, 120"
$ws.Range("D6").Value = "25, 30, 50, 80"
$ws.Range("D7").Value = "30, 50, 80"
$ws.Range("D8").Value = "40, 60, 80, 100"
$ws.Range("D9").Value = "40, 60, 80, 100"
$ws.Range("D10").Value = "60, 80, 100"
$ws.Range("D11").Value = "45, 75, 95, 115, 140"
$ws.Range("D12").Value = "20, 45, 65, 80, 90, 110"
$ws.Range("D13").Value = "70, 95, 120"
$ws.Range("D14").Value = "45, 70, 95, 120, 145, 195, 245"
$ws.Range("D15").Value = "25, 35, 50, 75, 100"

# Update comment text for the PhenolicFoam row (now column F)
$ws.Range("F12").Value = "two layers (K118)"

# Match the saved selection from the edited workbook
$ws.Range("D6").Select()
